# "Generate Report for Handoff" - refreshes the handoff/handback timestamps
# that were recomputed for the files currently "Ready for handoff" (and the
# one "Handback transform failed" entry), collapsing them onto the new
# uniform timestamps produced by this report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest Handoff Date" column (D) for rows 7, 10-16 ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-23-14 06:23:53"
$overview.Range("D10").Value = "2016-23-14 06:23:53"
$overview.Range("D11").Value = "2016-23-14 06:23:53"
$overview.Range("D12").Value = "2016-23-14 06:23:53"
$overview.Range("D13").Value = "2016-23-14 06:23:53"
$overview.Range("D14").Value = "2016-23-14 06:23:53"
$overview.Range("D15").Value = "2016-23-14 06:23:53"
$overview.Range("D16").Value = "2016-23-14 06:23:53"

# --- zh-cn sheet: "Latest Handoff Datetime" column (E) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-14 06:23:50"
$zhcn.Range("E10").Value = "2016-03-14 06:23:50"
$zhcn.Range("E11").Value = "2016-03-14 06:23:50"
$zhcn.Range("E12").Value = "2016-03-14 06:23:50"
$zhcn.Range("E13").Value = "2016-03-14 06:23:50"
$zhcn.Range("E14").Value = "2016-03-14 06:23:50"
$zhcn.Range("E15").Value = "2016-03-14 06:23:50"
$zhcn.Range("E16").Value = "2016-03-14 06:23:50"

# --- de-de sheet: "Latest Handoff Datetime" column (E) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-14 06:23:53"
$dede.Range("E10").Value = "2016-03-14 06:23:53"
$dede.Range("E11").Value = "2016-03-14 06:23:53"
$dede.Range("E12").Value = "2016-03-14 06:23:53"
$dede.Range("E13").Value = "2016-03-14 06:23:53"
$dede.Range("E14").Value = "2016-03-14 06:23:53"
$dede.Range("E15").Value = "2016-03-14 06:23:53"
$dede.Range("E16").Value = "2016-03-14 06:23:53"
